# Updated cryptos list on Sat Apr  1 19:49:05 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for each
# coin row (rows 2-51) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry holds the new Price (D) and Volume(1h) (E) text for a row.
# D (or E) is an empty string when that particular cell did not change.
$updates = @(
    @{ Row = 2; D = '28.413.68'; E = '  -0.11%  ' },
    @{ Row = 3; D = '1.815.67'; E = '  -0.66%  ' },
    @{ Row = 4; D = '1.003'; E = '  +0.24%  ' },
    @{ Row = 5; D = '314.88'; E = '  -0.74%  ' },
    @{ Row = 6; D = '1.001'; E = '  +0.17%  ' },
    @{ Row = 7; D = '0.5076'; E = '  -4.96%  ' },
    @{ Row = 8; D = '0.3950'; E = '  -1.71%  ' },
    @{ Row = 9; D = '0.08221'; E = '  +8.08%  ' },
    @{ Row = 10; D = '41.62'; E = '  -0.43%  ' },
    @{ Row = 11; D = '1.105'; E = '  -0.50%  ' },
    @{ Row = 12; D = '20.95'; E = '  -0.02%  ' },
    @{ Row = 13; D = '6.275'; E = '  -0.87%  ' },
    @{ Row = 14; D = ''; E = '  +0.18%  ' },
    @{ Row = 15; D = '7.491'; E = '  -1.68%  ' },
    @{ Row = 16; D = '1.813.76'; E = '  -1.06%  ' },
    @{ Row = 17; D = '0.00001141'; E = '  +6.23%  ' },
    @{ Row = 18; D = '92.33'; E = '  +3.08%  ' },
    @{ Row = 19; D = '0.06624'; E = '  +0.56%  ' },
    @{ Row = 20; D = '17.67'; E = '  -0.24%  ' },
    @{ Row = 21; D = ''; E = '  +0.14%  ' },
    @{ Row = 22; D = '6.093'; E = '  +0.16%  ' },
    @{ Row = 23; D = '28.445.93'; E = '  -0.04%  ' },
    @{ Row = 24; D = '11.29'; E = '  +0.96%  ' },
    @{ Row = 25; D = '2.268'; E = '  +2.46%  ' },
    @{ Row = 26; D = '21.15'; E = '  +2.34%  ' },
    @{ Row = 27; D = '155.33'; E = '  -1.25%  ' },
    @{ Row = 28; D = '2.025.42'; E = '  -0.87%  ' },
    @{ Row = 29; D = '2.398'; E = '  -2.33%  ' },
    @{ Row = 30; D = '125.82'; E = '  +1.42%  ' },
    @{ Row = 31; D = '1.107'; E = '  -1.43%  ' },
    @{ Row = 32; D = '0.1095'; E = '  -0.58%  ' },
    @{ Row = 33; D = '5.789'; E = '  +2.13%  ' },
    @{ Row = 34; D = '3.650'; E = '  +0.21%  ' },
    @{ Row = 35; D = '0.07028'; E = '  -6.00%  ' },
    @{ Row = 36; D = '0.2217'; E = '  -0.47%  ' },
    @{ Row = 37; D = '0.02335'; E = '  -0.40%  ' },
    @{ Row = 38; D = '5.205'; E = '  -0.05%  ' },
    @{ Row = 39; D = '8.823'; E = '  -0.22%  ' },
    @{ Row = 40; D = '0.6267'; E = '' },
    @{ Row = 41; D = '11.24'; E = '  -0.51%  ' },
    @{ Row = 42; D = '1.174'; E = '  -0.21%  ' },
    @{ Row = 43; D = ''; E = '  +0.15%  ' },
    @{ Row = 44; D = '1.403'; E = '  +0.84%  ' },
    @{ Row = 45; D = '13.49'; E = '  -0.29%  ' },
    @{ Row = 46; D = '3.740'; E = '  +1.00%  ' },
    @{ Row = 47; D = '0.5890'; E = '  +0.78%  ' },
    @{ Row = 48; D = '124.89'; E = '  +0.06%  ' },
    @{ Row = 49; D = '1.973'; E = '  -1.39%  ' },
    @{ Row = 50; D = '1.186'; E = '  -1.29%  ' },
    @{ Row = 51; D = '0.06885'; E = '  -0.06%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.D -ne "") {
        $dCell = $ws.Cells.Item($row, 4)
        # Values such as "1.003" or "314.88" look like plain numbers to Excel
        # and would otherwise be auto-converted from text into a numeric
        # value, dropping things like trailing zeros. Force the cell to text
        # first so the literal string is preserved, just like the original
        # formatted price strings (e.g. "28.427.41") already are.
        if ($u.D -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }

    if ($u.E -ne "") {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
